# "added calculator for upstream LCI CO2 emissions to factory"
#
# The underlying data edit (once the surrounding Excel-version metadata
# churn - Windows Excel -> Mac Excel re-save - is stripped away) is:
#   1. Rename the worksheet "Sheet1" -> "upstream".
#   2. Simplify three column headers on row 1:
#        B1: "upstream CO2"     -> "CO2"
#        C1: "upstream\nCH4"    -> "CH4"
#        E1: "CO2 removal"      -> "CO2 removed"
#   3. Rename the substance in A23:
#        A23: "iron ore - 65% Fe" -> "iron ore"
#
# All other hunks in the diff (fileVersion/xr namespaces/calcId/absPath/
# row-height & column-width rounding/dyDescent) are artifacts of the file
# having been re-saved by a different Excel build and aren't reachable
# (or meaningful) via the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "upstream"

$ws.Range("B1").Value = "CO2"
$ws.Range("C1").Value = "CH4"
$ws.Range("E1").Value = "CO2 removed"

$ws.Range("A23").Value = "iron ore"
